# Insert a new data row at row 21 (pushing existing rows 21:125 down to 22:126)
# and populate it with a new weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(21).Insert()

$ws.Range("A21").Value = 7
$ws.Range("B21").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C21").Value = "Ñuble"
$ws.Range("D21").Value = 44453
$ws.Range("E21").Value = 16
$ws.Range("F21").Value = 100112032
$ws.Range("G21").Value = "Zapallo italiano"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 120
$ws.Range("K21").Value = 16000
$ws.Range("L21").Value = 17000
$ws.Range("M21").Value = 16500
$ws.Range("N21").Value = "$/caja 50 unidades"
$ws.Range("O21").Value = "Región de Arica y Parinacota"
$ws.Range("P21").Value = 330
$ws.Range("Q21").Value = 50
$ws.Range("R21").Value = "Hortaliza"
